$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B2:E13 with unrounded floating point values (removing internal round/int casts)
$ws.Range("B2").Value = 591.1889053626218
$ws.Range("C2").Value = 370.7240588753156
$ws.Range("D2").Value = 312.0760440438579
$ws.Range("E2").Value = 274.6007624900782

$ws.Range("B3").Value = 706.539529825127
$ws.Range("C3").Value = 440.5443958823474
$ws.Range("D3").Value = 367.5404322105073
$ws.Range("E3").Value = 326.0277131780947

$ws.Range("B4").Value = 715.1048550385578
$ws.Range("C4").Value = 442.785858205543
$ws.Range("D4").Value = 367.1781300116164
$ws.Range("E4").Value = 325.1895444952996

$ws.Range("B5").Value = 544.295832905458
$ws.Range("C5").Value = 323.777547441535
$ws.Range("D5").Value = 257.8501257606252
$ws.Range("E5").Value = 228.2999937279845

$ws.Range("B6").Value = 481.1783049099772
$ws.Range("C6").Value = 285.498380909673
$ws.Range("D6").Value = 229.1516582824213
$ws.Range("E6").Value = 200.2069709732834

$ws.Range("B7").Value = 54.83687935742706
$ws.Range("C7").Value = 31.23958445848835
$ws.Range("D7").Value = 24.62755649136491
$ws.Range("E7").Value = 21.2165279986595

$ws.Range("B8").Value = 1331.242707437317
$ws.Range("C8").Value = 1096.965321119155
$ws.Range("D8").Value = 1036.678090764719
$ws.Range("E8").Value = 1058.730947649033

$ws.Range("B9").Value = 624.5632035571942
$ws.Range("C9").Value = 385.5847459911026
$ws.Range("D9").Value = 319.1188920320095
$ws.Range("E9").Value = 282.7667903645179

$ws.Range("B10").Value = 320.02988804888
$ws.Range("C10").Value = 179.2592538594663
$ws.Range("D10").Value = 142.6037805361524
$ws.Range("E10").Value = 121.2436951680541

$ws.Range("B11").Value = 60.50120700807447
$ws.Range("C11").Value = 31.59537102645282
$ws.Range("D11").Value = 24.69471916938721
$ws.Range("E11").Value = 22.62977408003405

$ws.Range("B12").Value = 139.3682673508567
$ws.Range("C12").Value = 82.3435814059286
$ws.Range("D12").Value = 65.23333897258398
$ws.Range("E12").Value = 54.80657424704891

$ws.Range("B13").Value = 179.2266158371581
$ws.Range("C13").Value = 99.82749828450733
$ws.Range("D13").Value = 80.20582558385682
$ws.Range("E13").Value = 69.00605023841192

$wb.Save()